$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New journal entry row 13: date, work description, hours spent.
$ws.Range("A13").Value = "2025-12-12"
$ws.Range("B13").Value = "Revue de la documentation mise en place le 5.12.2025"
$ws.Range("D13").Value = 1

# Reflect where the user's cursor ended up after typing the entry.
$ws.Range("B14:C14").Select()
